# Edit: split the run " con indirizzo di ritiro e indirizzo di consegna" into several
# runs that add "codice identificativo", "status" and "telefono Cliente e telefono Azienda"
# to the list of fields shown for a delivery (as described in the commit message).

$d = $word.ActiveDocument

$searchText = " con indirizzo di ritiro e indirizzo di consegna"

$findRange = $d.Content
$found = $findRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target text to edit: '$searchText'"
}

# Grab the whole paragraph that contains the text we need to change; InsertXML()
# always replaces the full paragraph it targets, so we rebuild that paragraph's
# content in full (unaffected runs are kept exactly as they were; only the one
# run that held the search text is split into the new, finer-grained runs).
$para = $findRange.Paragraphs(1)
$pRange = $para.Range

$expectedParaText = "Il Sistema mostra al Fattorino i dettagli della consegna con indirizzo di ritiro e indirizzo di consegna."
if ($pRange.Text.TrimEnd([char]7, [char]13, [char]10) -ne $expectedParaText) {
    throw ("Unexpected paragraph text, aborting: [" + $pRange.Text + "]")
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="22F35CB2" w14:textId="62C94D2B" w:rsidR="009D5A78" w:rsidRPr="009D5A78" w:rsidRDefault="009D5A78" w:rsidP="003446EB"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="374"/></w:pPr><w:r w:rsidRPr="009D5A78"><w:t xml:space="preserve">Il Sistema mostra al Fattorino </w:t></w:r><w:r w:rsidR="002D1157"><w:t xml:space="preserve">i dettagli </w:t></w:r><w:r w:rsidRPr="009D5A78"><w:t>dell</w:t></w:r><w:r w:rsidR="002D1157"><w:t>a</w:t></w:r><w:r w:rsidRPr="009D5A78"><w:t xml:space="preserve"> consegn</w:t></w:r><w:r w:rsidR="002D1157"><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> con indirizzo di ritiro</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> indirizzo di consegna</w:t></w:r><w:r><w:t>, codice identificativo</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> status</w:t></w:r><w:r><w:t>, telefono Cliente e telefono Azienda</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="002D1157"><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pRange.InsertXML($xml)
